$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 11 de Abril de 2020 a las 10:52'

# Row 19
$ws.Cells.Item(19, 4).Value = 6604
$ws.Cells.Item(19, 5).Value = 6726
$ws.Cells.Item(19, 6).Value = 246
$ws.Cells.Item(19, 7).Value = 18
$ws.Cells.Item(19, 8).Value = 337

# Row 21
$ws.Cells.Item(21, 1).Value = 'Israel'
$ws.Cells.Item(21, 2).Value = 10505
$ws.Cells.Item(21, 3).Value = 97
$ws.Cells.Item(21, 4).Value = 1236
$ws.Cells.Item(21, 5).Value = 9174
$ws.Cells.Item(21, 6).Value = 191
$ws.Cells.Item(21, 8).Value = 95

# Row 22
$ws.Cells.Item(22, 1).Value = 'Corea del Sur'
$ws.Cells.Item(22, 2).Value = 10480
$ws.Cells.Item(22, 3).Value = 30
$ws.Cells.Item(22, 4).Value = 7243
$ws.Cells.Item(22, 5).Value = 3026
$ws.Cells.Item(22, 6).Value = 55
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(22, 8).Value = 211

# Row 40
$ws.Cells.Item(40, 1).Value = 'Indonesia'
$ws.Cells.Item(40, 2).Value = 3842
$ws.Cells.Item(40, 3).Value = 330
$ws.Cells.Item(40, 4).Value = 286
$ws.Cells.Item(40, 5).Value = 3229
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 21
$ws.Cells.Item(40, 8).Value = 327

# Row 41
$ws.Cells.Item(41, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(41, 2).Value = 3651
$ws.Cells.Item(41, 4).Value = 685
$ws.Cells.Item(41, 5).Value = 2919
$ws.Cells.Item(41, 6).Value = 57
$ws.Cells.Item(41, 8).Value = 47

# Row 69
$ws.Cells.Item(69, 1).Value = 'Hong Kong'
$ws.Cells.Item(69, 2).Value = 1001
$ws.Cells.Item(69, 3).Value = 11
$ws.Cells.Item(69, 4).Value = 336
$ws.Cells.Item(69, 5).Value = 661
$ws.Cells.Item(69, 6).Value = 14
$ws.Cells.Item(69, 8).Value = 4

# Row 70
$ws.Cells.Item(70, 1).Value = 'Barein'
$ws.Cells.Item(70, 2).Value = 998
$ws.Cells.Item(70, 3).Value = 73
$ws.Cells.Item(70, 4).Value = 551
$ws.Cells.Item(70, 5).Value = 441
$ws.Cells.Item(70, 6).Value = 3
$ws.Cells.Item(70, 8).Value = 6

# Row 71
$ws.Cells.Item(71, 1).Value = 'Kuwait'
$ws.Cells.Item(71, 2).Value = 993
$ws.Cells.Item(71, 4).Value = 133
$ws.Cells.Item(71, 5).Value = 859
$ws.Cells.Item(71, 6).Value = 26
$ws.Cells.Item(71, 8).Value = 1

# Row 72
$ws.Cells.Item(72, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(72, 2).Value = 991
$ws.Cells.Item(72, 4).Value = 159
$ws.Cells.Item(72, 5).Value = 822
$ws.Cells.Item(72, 6).Value = 27
$ws.Cells.Item(72, 8).Value = 10

# Row 75
$ws.Cells.Item(75, 4).Value = 68
$ws.Cells.Item(75, 5).Value = 781

# Row 80
$ws.Cells.Item(80, 1).Value = 'Uzbekistan'
$ws.Cells.Item(80, 2).Value = 694
$ws.Cells.Item(80, 3).Value = 70
$ws.Cells.Item(80, 4).Value = 42
$ws.Cells.Item(80, 5).Value = 649
$ws.Cells.Item(80, 6).Value = 8
$ws.Cells.Item(80, 8).Value = 3

# Row 81
$ws.Cells.Item(81, 1).Value = 'Tunez'
$ws.Cells.Item(81, 2).Value = 671
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = 43
$ws.Cells.Item(81, 5).Value = 603
$ws.Cells.Item(81, 6).Value = 85
$ws.Cells.Item(81, 8).Value = 25

# Row 92
$ws.Cells.Item(92, 1).Value = 'Banglades'
$ws.Cells.Item(92, 2).Value = 482
$ws.Cells.Item(92, 3).Value = 58
$ws.Cells.Item(92, 4).Value = 36
$ws.Cells.Item(92, 5).Value = 416
$ws.Cells.Item(92, 6).Value = 1
$ws.Cells.Item(92, 7).Value = 3
$ws.Cells.Item(92, 8).Value = 30

# Row 93
$ws.Cells.Item(93, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(93, 2).Value = 480
$ws.Cells.Item(93, 4).Value = 54
$ws.Cells.Item(93, 5).Value = 423
$ws.Cells.Item(93, 8).Value = 3

# Row 94
$ws.Cells.Item(94, 1).Value = 'Burkina Faso'
$ws.Cells.Item(94, 2).Value = 448
$ws.Cells.Item(94, 4).Value = 149
$ws.Cells.Item(94, 5).Value = 273
$ws.Cells.Item(94, 8).Value = 26

# Row 95
$ws.Cells.Item(95, 1).Value = 'Niger'
$ws.Cells.Item(95, 2).Value = 438
$ws.Cells.Item(95, 4).Value = 41
$ws.Cells.Item(95, 5).Value = 386
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 8).Value = 11

# Row 170
$ws.Cells.Item(170, 1).Value = 'Laos'
$ws.Cells.Item(170, 3).Value = 2
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 18

# Row 171
$ws.Cells.Item(171, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(171, 4).Value = 1
$ws.Cells.Item(171, 5).Value = 17

# Row 172
$ws.Cells.Item(172, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(172, 2).Value = 18
$ws.Cells.Item(172, 4).Value = 3
$ws.Cells.Item(172, 5).Value = 15

# Row 173
$ws.Cells.Item(173, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 17
$ws.Cells.Item(173, 8).Value = 0

# Row 174
$ws.Cells.Item(174, 1).Value = 'Sudan'
$ws.Cells.Item(174, 2).Value = 17
$ws.Cells.Item(174, 4).Value = 2
$ws.Cells.Item(174, 5).Value = 13
$ws.Cells.Item(174, 8).Value = 2

# Row 175
$ws.Cells.Item(175, 1).Value = 'Fiyi'
